# Update cryptos list prices / 1h-volume percentages (GitHub Actions scrape refresh)
# Values are forced to text via a leading quote (matches source data stored as
# inline strings) then the style is reset to Normal so no stray number-format is
# left behind on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'59.229.01"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.54%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.526.76"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +1.34%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.04%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'534.70"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.16%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'139.78"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -3.17%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  +0.11%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.562"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -2.14%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'2.538.79"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +0.51%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.0990"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -0.91%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.160"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +1.74%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  -1.37%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.354"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +0.27%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'2.978.39"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +1.59%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'23.14"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -2.91%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'59.205.58"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +0.65%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'0.0000140"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +1.06%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'2.535.18"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +0.58%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'10.97"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -2.78%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'4.22"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -1.58%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'320.89"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -0.69%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'  +0.38%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'5.80"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +0.96%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'62.74"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +1.91%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.418"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -4.84%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  +2.63%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'0.996"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +0.17%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'7.82"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +0.82%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'6.74"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -1.02%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'0.0₃0769"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -0.88%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'1.79"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +0.22%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'161.19"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +0.88%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'0.998"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +0.29%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  -7.45%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'1.44"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -0.71%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'18.47"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -0.23%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  -4.64%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'1.58"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -2.46%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  +0.72%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'3.64"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -0.65%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'5.33"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -6.32%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'286.71"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -7.29%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.803"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -2.25%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.999"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +0.06%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.601"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +1.06%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'10.86"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +0.87%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'123.99"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -0.75%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.0927"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -0.38%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'18.60"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +0.25%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.0509"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -1.16%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.0222"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -2.59%  "
$ws.Range("E51").Style = "Normal"
